$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 393032.62
$ws.Range("J17").Value = 400873.3
$ws.Range("L17").Value = 1202619.9
$ws.Range("N17").Value = -1202955.9
# Row 46
$ws.Range("H46").Value = 600
$ws.Range("J46").Value = 600
$ws.Range("L46").Value = 1800
$ws.Range("N46").Value = -2038
# Row 48
$ws.Range("H48").Value = 6198
$ws.Range("I48").Value = 995
$ws.Range("J48").Value = 9666.667
$ws.Range("K48").Value = 2985
$ws.Range("L48").Value = 29000.001
$ws.Range("M48").Value = -2693
$ws.Range("N48").Value = -29584.001
# Row 55
$ws.Range("H55").Value = 169.2
$ws.Range("I55").Value = 86
$ws.Range("J55").Value = 502
$ws.Range("K55").Value = 86
$ws.Range("L55").Value = 502
$ws.Range("M55").Value = 128
$ws.Range("N55").Value = -930
# Row 56
$ws.Range("H56").Value = 6198
$ws.Range("I56").Value = 995
$ws.Range("J56").Value = 9666.667
$ws.Range("K56").Value = 2985
$ws.Range("L56").Value = 29000.001
$ws.Range("M56").Value = -2451
$ws.Range("N56").Value = -30068.001
# Row 58
$ws.Range("H58").Value = 137.75
$ws.Range("I58").Value = 137.75
$ws.Range("K58").Value = 413.25
$ws.Range("M58").Value = -263.25
# Row 60
$ws.Range("H60").Value = 600
$ws.Range("J60").Value = 600
$ws.Range("L60").Value = 1800
$ws.Range("N60").Value = -2768
# Row 61
$ws.Range("H61").Value = 16933.334
$ws.Range("I61").Value = 16933.334
$ws.Range("K61").Value = 50800.00199999999
$ws.Range("M61").Value = -50628.00199999999
# Row 64
$ws.Range("H64").Value = 111115110
$ws.Range("I64").Value = 111115110
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 111115110
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -111114862
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 111115110
$ws.Range("I67").Value = 111115110
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 111115110
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -111114252
$ws.Range("N67").ClearContents()
# Row 70
$ws.Range("H70").Value = 2911.1365
$ws.Range("I70").Value = 2258
$ws.Range("J70").Value = 3156.0625
$ws.Range("K70").Value = 6774
$ws.Range("L70").Value = 9468.1875
$ws.Range("M70").Value = -6504
$ws.Range("N70").Value = -10008.1875
# Row 73
$ws.Range("H73").Value = 2911.1365
$ws.Range("I73").Value = 2258
$ws.Range("J73").Value = 3156.0625
$ws.Range("K73").Value = 6774
$ws.Range("L73").Value = 9468.1875
$ws.Range("M73").Value = -5838
$ws.Range("N73").Value = -11340.1875

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2059526.2
$ws.Range("I32").Value = 941.6087
$ws.Range("J32").Value = 13896388
$ws.Range("K32").Value = 941.6087
$ws.Range("L32").Value = 13896388
$ws.Range("M32").Value = -654.6087
$ws.Range("N32").Value = -13896962
# Row 63
$ws.Range("H63").Value = 3324.4
$ws.Range("I63").Value = 2207.3333
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2207.3333
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1521.3333
$ws.Range("N63").Value = -6372
# Row 66
$ws.Range("H66").Value = 3324.4
$ws.Range("I66").Value = 2207.3333
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 11036.6665
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -7604.666499999999
$ws.Range("N66").Value = -31864
# Row 110
$ws.Range("H110").Value = 6512.8887
$ws.Range("I110").Value = 3766
$ws.Range("J110").Value = 7062.2666
$ws.Range("K110").Value = 3766
$ws.Range("L110").Value = 7062.2666
$ws.Range("M110").Value = -1721
$ws.Range("N110").Value = -11152.2666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4214
$ws.Range("I105").Value = 3126.4167
$ws.Range("J105").Value = 6078.4287
$ws.Range("K105").Value = 3126.4167
$ws.Range("L105").Value = 6078.4287
$ws.Range("M105").Value = -1379.4167
$ws.Range("N105").Value = -9572.4287

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 167445.56
$ws.Range("I5").Value = 583
$ws.Range("J5").Value = 358145.66
$ws.Range("K5").Value = 1749
$ws.Range("L5").Value = 1074436.98
$ws.Range("M5").Value = -1637
$ws.Range("N5").Value = -1074660.98
# Row 38
$ws.Range("H38").Value = 2090.7144
$ws.Range("J38").Value = 2061
$ws.Range("L38").Value = 6183
$ws.Range("N38").Value = -6877
# Row 135
$ws.Range("H135").Value = 167445.56
$ws.Range("I135").Value = 583
$ws.Range("J135").Value = 358145.66
$ws.Range("K135").Value = 5247
$ws.Range("L135").Value = 3223310.94
$ws.Range("M135").Value = -2712
$ws.Range("N135").Value = -3228380.94
# Row 137
$ws.Range("H137").Value = 1935.1111
$ws.Range("J137").Value = 1912.4286
$ws.Range("L137").Value = 5737.2858
$ws.Range("N137").Value = -15937.2858

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Range("H49").Value = 15861
$ws.Range("J49").Value = 15861
$ws.Range("L49").Value = 15861
$ws.Range("N49").Value = -16229
# Row 80
$ws.Range("H80").Value = 4686.1904
$ws.Range("I80").Value = 3626.7856
$ws.Range("J80").Value = 6805
$ws.Range("K80").Value = 3626.7856
$ws.Range("L80").Value = 6805
$ws.Range("M80").Value = -2628.7856
$ws.Range("N80").Value = -8801
# Row 83
$ws.Range("H83").Value = 4686.1904
$ws.Range("I83").Value = 3626.7856
$ws.Range("J83").Value = 6805
$ws.Range("K83").Value = 18133.928
$ws.Range("L83").Value = 34025
$ws.Range("M83").Value = -13141.928
$ws.Range("N83").Value = -44009

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 3625.6086
$ws.Range("I82").Value = 1085.1428
$ws.Range("J82").Value = 7577.4443
$ws.Range("K82").Value = 1085.1428
$ws.Range("L82").Value = 7577.4443
$ws.Range("M82").Value = -724.1428000000001
$ws.Range("N82").Value = -8299.4443
# Row 85
$ws.Range("H85").Value = 3625.6086
$ws.Range("I85").Value = 1085.1428
$ws.Range("J85").Value = 7577.4443
$ws.Range("K85").Value = 1085.1428
$ws.Range("L85").Value = 7577.4443
$ws.Range("M85").Value = 162.8571999999999
$ws.Range("N85").Value = -10073.4443
# Row 100
$ws.Range("H100").Value = 4676.3335
$ws.Range("I100").Value = 8021.75
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 8021.75
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -7480.75
$ws.Range("N100").Value = -3082
# Row 132
$ws.Range("H132").Value = 3227.0925
$ws.Range("I132").Value = 2257.6223
$ws.Range("J132").Value = 8074.4443
$ws.Range("K132").Value = 6772.8669
$ws.Range("L132").Value = 24223.3329
$ws.Range("M132").Value = -4242.8669
$ws.Range("N132").Value = -29283.3329

